# Weekly fruit/vegetable price update for "Higo" (Vega Central Mapocho de
# Santiago). The diff shows the data rows (4-15, grouped in pairs per
# quality "Primera"/"Segunda") being re-shuffled to different dates while
# keeping each pair's Volumen/Precio/Origen values bundled together.
#
# Target values per row (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg, Origen):

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = @{ D = 44322; M = 45;  N = 12000; O = 12000; P = 12000; S = 1714; R = "Región Metropolitana" }
    5  = @{ D = 44322; M = 80;  N = 8000;  O = 8000;  P = 8000;  S = 1143; R = "Región Metropolitana" }
    6  = @{ D = 44299; M = 80;  N = 15000; O = 15000; P = 15000; S = 2143; R = "Provincia de Santiago" }
    7  = @{ D = 44299; M = 75;  N = 12000; O = 12000; P = 12000; S = 1714; R = "Provincia de Santiago" }
    8  = @{ D = 44292; M = 25;  N = 16000; O = 16000; P = 16000; S = 2286; R = "Región Metropolitana" }
    9  = @{ D = 44292; M = 30;  N = 15000; O = 15000; P = 15000; S = 2143; R = "Región Metropolitana" }
    10 = @{ D = 44300; M = 100; N = 15000; O = 15000; P = 15000; S = 2143; R = "Región Metropolitana" }
    11 = @{ D = 44300; M = 80;  N = 12000; O = 12000; P = 12000; S = 1714; R = "Región Metropolitana" }
    12 = @{ D = 44302; M = 50;  N = 15000; O = 15000; P = 15000; S = 2143; R = "Región Metropolitana" }
    13 = @{ D = 44302; M = 30;  N = 12000; O = 12000; P = 12000; S = 1714; R = "Región Metropolitana" }
    14 = @{ D = 44301; M = 100; N = 14000; O = 14000; P = 14000; S = 2000; R = "Región Metropolitana" }
    15 = @{ D = 44301; M = 80;  N = 12000; O = 12000; P = 12000; S = 1714; R = "Región Metropolitana" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
